$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("D9").Value = 21168.54081889284
$ws.Range("D10").Value = 21168.54081889284
$ws.Range("D14").Value = 20217.59999999929
$ws.Range("D15").Value = 20217.59999999929
